$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 683.8182
$ws.Range("J32").Value = 660.375
$ws.Range("L32").Value = 660.375
$ws.Range("N32").Value = -1312.375
$ws.Range("H51").Value = 4659.8
$ws.Range("I51").Value = 4274.5
$ws.Range("J51").Value = 4916.6665
$ws.Range("K51").Value = 4274.5
$ws.Range("L51").Value = 4916.6665
$ws.Range("M51").Value = -3790.5
$ws.Range("N51").Value = -5884.6665
$ws.Range("H62").Value = 2724.6667
$ws.Range("I62").Value = 1994.75
$ws.Range("K62").Value = 1994.75
$ws.Range("M62").Value = -1370.75
$ws.Range("H65").Value = 2724.6667
$ws.Range("I65").Value = 1994.75
$ws.Range("K65").Value = 9973.75
$ws.Range("M65").Value = -6853.75
$ws.Range("H74").Value = 4431.875
$ws.Range("I74").Value = 4509.6665
$ws.Range("K74").Value = 4509.6665
$ws.Range("M74").Value = -3573.6665
$ws.Range("H77").Value = 4431.875
$ws.Range("I77").Value = 4509.6665
$ws.Range("K77").Value = 22548.3325
$ws.Range("M77").Value = -17868.3325
$ws.Range("H96").Value = 2850
$ws.Range("I96").Value = 750
$ws.Range("J96").Value = 4950
$ws.Range("K96").Value = 2250
$ws.Range("L96").Value = 14850
$ws.Range("M96").Value = -877
$ws.Range("N96").Value = -17596
$ws.Range("H113").Value = 11634.917
$ws.Range("I113").Value = 15935.75
$ws.Range("K113").Value = 15935.75
$ws.Range("M113").Value = -12681.75
$ws.Range("H116").Value = 12096.363
$ws.Range("I116").Value = 22394.8
$ws.Range("K116").Value = 22394.8
$ws.Range("M116").Value = -18952.8
$ws.Range("H117").Value = 30000
$ws.Range("J117").Value = 30000
$ws.Range("L117").Value = 30000
$ws.Range("N117").Value = -39178
$ws.Range("H131").Value = 1025.6923
$ws.Range("I131").Value = 527.8333
$ws.Range("J131").Value = 7000
$ws.Range("K131").Value = 1583.4999
$ws.Range("L131").Value = 21000
$ws.Range("M131").Value = 3456.5001
$ws.Range("N131").Value = -31080
$ws.Range("H135").Value = 675.6923
$ws.Range("I135").Value = 553.4
$ws.Range("K135").Value = 4980.599999999999
$ws.Range("M135").Value = -2445.599999999999
$ws.Range("H137").Value = 2087.1667
$ws.Range("I137").Value = 1982.8889
$ws.Range("K137").Value = 5948.6667
$ws.Range("M137").Value = -3398.6667
$ws.Range("H141").Value = 3645.4119
$ws.Range("I141").Value = 2637.2
$ws.Range("J141").Value = 5085.7144
$ws.Range("K141").Value = 7911.599999999999
$ws.Range("L141").Value = 15257.1432
$ws.Range("M141").Value = -2731.599999999999
$ws.Range("N141").Value = -25617.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4116.25
$ws.Range("I32").Value = 2906
$ws.Range("K32").Value = 2906
$ws.Range("M32").Value = -2619
$ws.Range("H45").Value = 4738237.5
$ws.Range("I45").Value = 6429872.5
$ws.Range("J45").Value = 1658.6
$ws.Range("K45").Value = 6429872.5
$ws.Range("L45").Value = 1658.6
$ws.Range("M45").Value = -6429495.5
$ws.Range("N45").Value = -2412.6
$ws.Range("H61").Value = 2901.25
$ws.Range("I61").Value = 2072.6667
$ws.Range("K61").Value = 2072.6667
$ws.Range("M61").Value = -1860.6667
$ws.Range("H86").Value = 25000
$ws.Range("I86").Value = 25000
$ws.Range("K86").Value = 25000
$ws.Range("H88").Value = 2801.5881
$ws.Range("J88").Value = 3550
$ws.Range("L88").Value = 3550
$ws.Range("N88").Value = -4362
$ws.Range("H89").Value = 25000
$ws.Range("I89").Value = 25000
$ws.Range("K89").Value = 75000
$ws.Range("H91").Value = 2801.5881
$ws.Range("J91").Value = 3550
$ws.Range("L91").Value = 3550
$ws.Range("N91").Value = -6358
$ws.Range("H122").Value = 1976.1428
$ws.Range("I122").Value = 1976.1428
$ws.Range("K122").Value = 5928.428400000001
$ws.Range("M122").Value = -3478.428400000001
$ws.Range("H136").Value = 2901.25
$ws.Range("I136").Value = 2072.6667
$ws.Range("K136").Value = 6218.000100000001
$ws.Range("M136").Value = -3668.000100000001
$ws.Range("M86").Value = -23814
$ws.Range("M89").Value = -69072

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 34250
$ws.Range("J87").Value = 28500
$ws.Range("L87").Value = 28500
$ws.Range("N87").Value = -30996
$ws.Range("H90").Value = 34250
$ws.Range("J90").Value = 28500
$ws.Range("L90").Value = 85500
$ws.Range("N90").Value = -97980

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3442.389
$ws.Range("I31").Value = 1367.5
$ws.Range("K31").Value = 1367.5
$ws.Range("M31").Value = -1072.5
$ws.Range("H34").Value = 3442.389
$ws.Range("I34").Value = 1367.5
$ws.Range("K34").Value = 1367.5
$ws.Range("M34").Value = -1165.5
$ws.Range("H99").Value = 1490.7
$ws.Range("I99").Value = 1378.5555
$ws.Range("K99").Value = 1378.5555
$ws.Range("M99").Value = 119.4445000000001
$ws.Range("H122").Value = 2093.8235
$ws.Range("I122").Value = 2800.4
$ws.Range("J122").Value = 1799.4166
$ws.Range("K122").Value = 8401.200000000001
$ws.Range("L122").Value = 5398.2498
$ws.Range("M122").Value = -5951.200000000001
$ws.Range("N122").Value = -10298.2498
$ws.Range("H126").Value = 1490.7
$ws.Range("I126").Value = 1378.5555
$ws.Range("K126").Value = 4135.666499999999
$ws.Range("M126").Value = -1665.666499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 515.2857
$ws.Range("J107").Value = 592.1818
$ws.Range("L107").Value = 1776.5454
$ws.Range("N107").Value = -5616.5454
$ws.Range("H131").Value = 11380534
$ws.Range("J131").Value = 20536.277
$ws.Range("L131").Value = 61608.83099999999
$ws.Range("N131").Value = -71688.83099999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4627.5
$ws.Range("J70").Value = 3999.25
$ws.Range("L70").Value = 3999.25
$ws.Range("N70").Value = -4539.25
$ws.Range("H73").Value = 4627.5
$ws.Range("J73").Value = 3999.25
$ws.Range("L73").Value = 3999.25
$ws.Range("N73").Value = -5871.25
$ws.Range("H102").Value = 1715.6296
$ws.Range("I102").Value = 1788
$ws.Range("K102").Value = 1788
$ws.Range("M102").Value = -166
$ws.Range("H132").Value = 3041.96
$ws.Range("J132").Value = 3256.6667
$ws.Range("L132").Value = 9770.000100000001
$ws.Range("N132").Value = -14830.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 528.9545000000001
$ws.Range("I107").Value = 449.26315
$ws.Range("K107").Value = 1347.78945
$ws.Range("M107").Value = 572.21055
$ws.Range("H122").Value = 53274.066
$ws.Range("I122").Value = 61241.69
$ws.Range("K122").Value = 183725.07
$ws.Range("M122").Value = -181275.07
